$d = $word.ActiveDocument

# Collapse to the very end of the document body content (right before sectPr).
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd

# Build a minimal OOXML package fragment describing the new paragraph:
#   "Tout commence quand…" followed by a _GoBack bookmark, with no extra
# direct formatting (matches what Word leaves behind after typing at the
# end of a document). Using InsertXML avoids inheriting the bold/underline/
# centered formatting of the preceding paragraph.
$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r><w:t>Tout commence quand&#8230;</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$endRange.InsertXML($newParagraphXml) | Out-Null
